{"js": "// Remove the trailing site-footer boilerplate that had been appended at the\n// end of the document (\"Ver no Jupiter Salvar em pdf Salvar em docx\" and the\n// \"\u00a9 2020 ...\" copyright line), along with the blank paragraph that\n// separated them from the preceding \"Requisitos\" text. The blank paragraph\n// that originally sat right after the \"LOM3016: ...\" line, as well as the\n// final page-break paragraph, are left untouched.\n\nconst body = context.document.body;\n\n// Locate the \"Ver no Jupiter ...\" paragraph via text search so the edit is\n// resilient to exact paragraph indices.\nconst results = body.search(\"Ver no Jupiter Salvar em pdf Salvar em docx\", {\n  matchCase: true\n});\nawait context.sync();\n\nif (results.items.length > 0) {\n  const jupiterPara = results.items[0].paragraphs.getFirst();\n  const blankPara = jupiterPara.getPrevious();\n  const copyrightPara = jupiterPara.getNext();\n  await context.sync();\n\n  // Delete the blank separator line plus the two footer paragraphs.\n  blankPara.delete();\n  jupiterPara.delete();\n  copyrightPara.delete();\n  await context.sync();\n}\n", "ps1": "# Remove the trailing site-footer boilerplate that had been appended at the\n# end of the document (\"Ver no Jupiter Salvar em pdf Salvar em docx\" and the\n# \"\u00a9 2020 ...\" copyright line), along with the blank paragraph that\n# separated them from the preceding \"Requisitos\" text. The blank paragraph\n# that originally sat right after the \"LOM3016: ...\" line, as well as the\n# final page-break paragraph, are left untouched.\n\n$d = $word.ActiveDocument\n\n# Locate the \"Ver no Jupiter ...\" paragraph by scanning paragraph text so the\n# edit is resilient to exact paragraph indices.\n$targetIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"Ver no Jupiter Salvar em pdf Salvar em docx*\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -gt 0) {\n    # Delete highest index first so earlier indices stay valid:\n    #   targetIndex + 1 -> the \"\u00a9 2020 ...\" copyright paragraph\n    #   targetIndex     -> the \"Ver no Jupiter ...\" paragraph\n    #   targetIndex - 1 -> the blank separator paragraph before it\n    $d.Paragraphs.Item($targetIndex + 1).Range.Delete()\n    $d.Paragraphs.Item($targetIndex).Range.Delete()\n    $d.Paragraphs.Item($targetIndex - 1).Range.Delete()\n}\n"}
